$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as plain text (inlineStr) in the source
# workbook, e.g. "246.13". Excel's Range.Value setter auto-coerces a
# numeric-looking string into a real number (and drops any significant
# trailing zeros), which would corrupt these text cells. Forcing the
# value in through Formula with a leading apostrophe (the same trick a
# user gets from typing a quote-prefixed entry) keeps it as literal
# text; resetting Style back to "Normal" afterwards strips the quote-
# prefix/Text-number-format styling that operation applies, so the
# cell ends up identical in type/format to how it started (plain text,
# default style) - just with the new value.
function Set-TextPrice($address, $value) {
    $c = $ws.Range($address)
    $c.Formula = "'" + $value
    $c.Style = "Normal"
}

# Row 2 (BNB)
Set-TextPrice "D2" "246.03"
# Row 3 (OKB)
Set-TextPrice "D3" "22.07"
# Row 4 (HuobiToken)
Set-TextPrice "D4" "5.371"
# Row 5 (Cronos)
Set-TextPrice "D5" "0.05864"
# Row 6
Set-TextPrice "D6" "3.385"
# Row 7
Set-TextPrice "D7" "6.387"
# Row 8
Set-TextPrice "D8" "0.8127"
# Row 9
Set-TextPrice "D9" "0.9618"
# Row 10
Set-TextPrice "D10" "0.1419"
# Row 11
Set-TextPrice "D11" "0.03582"
# Row 12
Set-TextPrice "D12" "0.07300"
# Row 14
Set-TextPrice "D14" "4.462"
# Row 15
Set-TextPrice "D15" "0.09393"
# Row 16
Set-TextPrice "D16" "0.001602"
# Row 17
Set-TextPrice "D17" "0.04838"
# Row 18 (One) - price + volume label change
Set-TextPrice "D18" "0.0005895"
$ws.Range("E18").Value = "17OneONEWorstin24h"
# Row 19
Set-TextPrice "D19" "0.006206"
# Row 20
Set-TextPrice "D20" "0.004084"
# Row 21
Set-TextPrice "D21" "0.0009819"
# Row 22
Set-TextPrice "D22" "0.00009704"
# Row 23
Set-TextPrice "D23" "3.691"
# Row 24
Set-TextPrice "D24" "2.198"
# Row 40 (IDEX)
Set-TextPrice "D40" "0.03851"

# Rows 41-43 got reshuffled (coin ranking reorder):
#   old 41 BKEXToken -> new 41 KickToken
#   old 42 CEJI      -> new 42 BKEXToken
#   old 43 KickToken -> new 43 CEJI
# Row 41 -> KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextPrice "D41" "0.006592"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 -> BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextPrice "D42" "0.1072"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 -> CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextPrice "D43" "0.003001"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 44
Set-TextPrice "D44" "0.005759"
# Row 45
Set-TextPrice "D45" "0.00005665"
# Row 47
Set-TextPrice "D47" "0.6516"
# Row 48
Set-TextPrice "D48" "0.07664"
# Row 50
Set-TextPrice "D50" "0.01010"
